$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sku values (A2:A8) with the new descriptive SKU codes ---
$ws.Range("A2").Value = "ST-TS-1-DAWG-BLK-L"
$ws.Range("A3").Value = "ST-TS-4-CRACK-GRN-M"
$ws.Range("A4").Value = "ST-TS-6-HEART-WHT-L"
$ws.Range("A5").Value = "ST-TS-1-DAWG-BLK-XL"
$ws.Range("A6").Value = "ST-TS-5-YAP-BLK-XL"
$ws.Range("A7").Value = "ST-TS-2-DAWG-GRN-XL"
$ws.Range("A8").Value = "ST-TS-3-CRACK-BLK-L"

# --- Add the new row (9) with its sku + qty ---
$ws.Range("A9").Value = "ST-TS-6-HEART-WHT-XL"
$ws.Range("B9").Value = 2

# --- Grow the table so the new row is included ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B10"))

# --- Highlight/center the qty for the HEART-WHT-L row ---
$ws.Range("B4").Interior.ColorIndex = -4142
$ws.Range("B4").HorizontalAlignment = -4108

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 22.140625
$ws.Columns("B").ColumnWidth = 8

# --- View: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 235
$ws.Range("D4").Select()
